# Adds three new sheets (CypherOutput_Message, StatOutput, StatOutput_Message)
# that report on a Neo4j "stats" Cypher query (file/sample/case/study counts
# for the Akita breed), mirroring the existing CypherOutput/Message pair that
# already reports the per-case Beagle query results.

$wb = $excel.ActiveWorkbook

$cypherOutput = $wb.Worksheets.Item("CypherOutput")
$message      = $wb.Worksheets.Item("Message")

# ---------------------------------------------------------------------
# 1) CypherOutput_Message - exact duplicate of the "Message" sheet
# ---------------------------------------------------------------------
$cypherOutputMessage = $wb.Worksheets.Add([System.Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$cypherOutputMessage.Name = "CypherOutput_Message"

for ($r = 1; $r -le 10; $r++) {
    $cypherOutputMessage.Cells.Item($r, 1).Value = $message.Cells.Item($r, 1).Value2
}

# ---------------------------------------------------------------------
# 2) StatOutput - header row + one row of counts, result of the stats query
# ---------------------------------------------------------------------
$statOutput = $wb.Worksheets.Add([System.Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$statOutput.Name = "StatOutput"

$statOutput.Cells.Item(1, 1).Value = "number_of_files"
$statOutput.Cells.Item(1, 2).Value = "number_of_sample"
$statOutput.Cells.Item(1, 3).Value = "number_of_cases"
$statOutput.Cells.Item(1, 4).Value = "number_of_study"

$statOutput.Cells.Item(2, 1).NumberFormat = "@"
$statOutput.Cells.Item(2, 1).Value = "1"
$statOutput.Cells.Item(2, 2).NumberFormat = "@"
$statOutput.Cells.Item(2, 2).Value = "2"
$statOutput.Cells.Item(2, 3).NumberFormat = "@"
$statOutput.Cells.Item(2, 3).Value = "1"
$statOutput.Cells.Item(2, 4).NumberFormat = "@"
$statOutput.Cells.Item(2, 4).Value = "1"

# ---------------------------------------------------------------------
# 3) StatOutput_Message - Message sheet content repeated twice; the second
#    copy's Cypher cell (row 18) holds the stats-query text (Akita breed)
#    instead of the per-case query text used by CypherOutput/row 8.
# ---------------------------------------------------------------------
$statOutputMessage = $wb.Worksheets.Add([System.Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$statOutputMessage.Name = "StatOutput_Message"

$statsCypher = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Akita']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

for ($r = 1; $r -le 10; $r++) {
    $statOutputMessage.Cells.Item($r, 1).Value = $message.Cells.Item($r, 1).Value2
}
for ($r = 1; $r -le 7; $r++) {
    $statOutputMessage.Cells.Item(10 + $r, 1).Value = $message.Cells.Item($r, 1).Value2
}
$statOutputMessage.Cells.Item(18, 1).Value = $statsCypher
$statOutputMessage.Cells.Item(19, 1).Value = $message.Cells.Item(9, 1).Value2
$statOutputMessage.Cells.Item(20, 1).Value = $message.Cells.Item(10, 1).Value2

# ---------------------------------------------------------------------
# Restore original active sheet/selection
# ---------------------------------------------------------------------
$cypherOutput.Activate()
$cypherOutput.Range("A1").Select()
